# Fruta / hortaliza, semanal
# Update the weekly price/date data in rows 2-18 (skipping unchanged rows 6 and 7)
# Columns: D = Fecha (date serial), J = Volumen, K = Precio minimo,
#          L = Precio maximo, M = Precio promedio ponderado, P = Precio $/Kg

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = 44893; J = 3300; K = 1200; L = 1300; M = 1261; P = 1261 }
    3  = @{ D = 45210; J = 550;  K = 1500; L = 1600; M = 1536; P = 1536 }
    4  = @{ D = 45203; J = 800;  K = 1800; L = 2000; M = 1900; P = 1900 }
    5  = @{ D = 45205; J = 3500; K = 1400; L = 1500; M = 1457; P = 1457 }
    8  = @{ D = 45062; J = 1700; K = 2800; L = 3000; M = 2900; P = 2900 }
    9  = @{ D = 44895; J = 200;  K = 1200; L = 1300; M = 1255; P = 1255 }
    10 = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
    11 = @{ D = 45212; J = 750;  K = 1400; L = 1500; M = 1440; P = 1440 }
    12 = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    13 = @{ D = 45233; J = 1050; K = 1400; L = 1500; M = 1438; P = 1438 }
    14 = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
    15 = @{ D = 45132; J = 170;  K = 2200; L = 2500; M = 2359; P = 2359 }
    16 = @{ D = 44907; J = 2300; K = 900;  L = 1000; M = 952;  P = 952 }
    17 = @{ D = 44883; J = 290;  K = 1400; L = 1500; M = 1434; P = 1434 }
    18 = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
